$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "acheron"
$ws.Cells.Item(2, 2).Value = "Nihility"
$ws.Cells.Item(2, 3).Value = "5"
$ws.Cells.Item(2, 4).Value = "Lightning"
$ws.Cells.Item(3, 1).Value = "argenti"
$ws.Cells.Item(3, 2).Value = "Erudition"
$ws.Cells.Item(3, 3).Value = "5"
$ws.Cells.Item(3, 4).Value = "Physical"
$ws.Cells.Item(4, 1).Value = "arlan"
$ws.Cells.Item(4, 2).Value = "Destruction"
$ws.Cells.Item(4, 3).Value = "4"
$ws.Cells.Item(4, 4).Value = "Lightning"
$ws.Cells.Item(5, 1).Value = "asta"
$ws.Cells.Item(5, 2).Value = "Harmony"
$ws.Cells.Item(5, 3).Value = "4"
$ws.Cells.Item(5, 4).Value = "Fire"
$ws.Cells.Item(6, 1).Value = "aventurine"
$ws.Cells.Item(6, 2).Value = "Preservation"
$ws.Cells.Item(6, 3).Value = "5"
$ws.Cells.Item(6, 4).Value = "Imaginary"
$ws.Cells.Item(7, 1).Value = "bailu"
$ws.Cells.Item(7, 2).Value = "Abundance"
$ws.Cells.Item(7, 3).Value = "5"
$ws.Cells.Item(7, 4).Value = "Lightning"
$ws.Cells.Item(8, 1).Value = "black-swan"
$ws.Cells.Item(8, 2).Value = "Nihility"
$ws.Cells.Item(8, 3).Value = "5"
$ws.Cells.Item(8, 4).Value = "Wind"
$ws.Cells.Item(9, 1).Value = "blade"
$ws.Cells.Item(9, 2).Value = "Destruction"
$ws.Cells.Item(9, 3).Value = "5"
$ws.Cells.Item(9, 4).Value = "Wind"
$ws.Cells.Item(10, 1).Value = "bronya"
$ws.Cells.Item(10, 2).Value = "Harmony"
$ws.Cells.Item(10, 3).Value = "5"
$ws.Cells.Item(10, 4).Value = "Wind"
$ws.Cells.Item(11, 1).Value = "clara"
$ws.Cells.Item(11, 2).Value = "Destruction"
$ws.Cells.Item(11, 3).Value = "5"
$ws.Cells.Item(11, 4).Value = "Physical"
$ws.Cells.Item(12, 1).Value = "dan-heng"
$ws.Cells.Item(12, 2).Value = "Hunt"
$ws.Cells.Item(12, 3).Value = "4"
$ws.Cells.Item(12, 4).Value = "Wind"
$ws.Cells.Item(13, 1).Value = "imbibitor-lunae"
$ws.Cells.Item(13, 2).Value = "Destruction"
$ws.Cells.Item(13, 3).Value = "5"
$ws.Cells.Item(13, 4).Value = "Imaginary"
$ws.Cells.Item(14, 1).Value = "dr-ratio"
$ws.Cells.Item(14, 2).Value = "Hunt"
$ws.Cells.Item(14, 3).Value = "5"
$ws.Cells.Item(14, 4).Value = "Imaginary"
$ws.Cells.Item(15, 1).Value = "fu-xuan"
$ws.Cells.Item(15, 2).Value = "Preservation"
$ws.Cells.Item(15, 3).Value = "5"
$ws.Cells.Item(15, 4).Value = "Quantum"
$ws.Cells.Item(16, 1).Value = "gallagher"
$ws.Cells.Item(16, 2).Value = "Abundance"
$ws.Cells.Item(16, 3).Value = "4"
$ws.Cells.Item(16, 4).Value = "Fire"
$ws.Cells.Item(17, 1).Value = "gepard"
$ws.Cells.Item(17, 2).Value = "Preservation"
$ws.Cells.Item(17, 3).Value = "5"
$ws.Cells.Item(17, 4).Value = "Ice"
$ws.Cells.Item(18, 1).Value = "guinaifen"
$ws.Cells.Item(18, 2).Value = "Nihility"
$ws.Cells.Item(18, 3).Value = "4"
$ws.Cells.Item(18, 4).Value = "Fire"
$ws.Cells.Item(19, 1).Value = "hanya"
$ws.Cells.Item(19, 2).Value = "Harmony"
$ws.Cells.Item(19, 3).Value = "4"
$ws.Cells.Item(19, 4).Value = "Physical"
$ws.Cells.Item(20, 1).Value = "herta"
$ws.Cells.Item(20, 2).Value = "Erudition"
$ws.Cells.Item(20, 3).Value = "4"
$ws.Cells.Item(20, 4).Value = "Ice"
$ws.Cells.Item(21, 1).Value = "himeko"
$ws.Cells.Item(21, 2).Value = "Erudition"
$ws.Cells.Item(21, 3).Value = "5"
$ws.Cells.Item(21, 4).Value = "Fire"
$ws.Cells.Item(22, 1).Value = "hook"
$ws.Cells.Item(22, 2).Value = "Destruction"
$ws.Cells.Item(22, 3).Value = "4"
$ws.Cells.Item(22, 4).Value = "Fire"
$ws.Cells.Item(23, 1).Value = "huohuo"
$ws.Cells.Item(23, 2).Value = "Abundance"
$ws.Cells.Item(23, 3).Value = "5"
$ws.Cells.Item(23, 4).Value = "Wind"
$ws.Cells.Item(24, 1).Value = "jing-yuan"
$ws.Cells.Item(24, 2).Value = "Erudition"
$ws.Cells.Item(24, 3).Value = "5"
$ws.Cells.Item(24, 4).Value = "Lightning"
$ws.Cells.Item(25, 1).Value = "jingliu"
$ws.Cells.Item(25, 2).Value = "Destruction"
$ws.Cells.Item(25, 3).Value = "5"
$ws.Cells.Item(25, 4).Value = "Ice"
$ws.Cells.Item(26, 1).Value = "kafka"
$ws.Cells.Item(26, 2).Value = "Nihility"
$ws.Cells.Item(26, 3).Value = "5"
$ws.Cells.Item(26, 4).Value = "Lightning"
$ws.Cells.Item(27, 1).Value = "luka"
$ws.Cells.Item(27, 2).Value = "Nihility"
$ws.Cells.Item(27, 3).Value = "4"
$ws.Cells.Item(27, 4).Value = "Physical"
$ws.Cells.Item(28, 1).Value = "luocha"
$ws.Cells.Item(28, 2).Value = "Abundance"
$ws.Cells.Item(28, 3).Value = "5"
$ws.Cells.Item(28, 4).Value = "Imaginary"
$ws.Cells.Item(29, 1).Value = "lynx"
$ws.Cells.Item(29, 2).Value = "Abundance"
$ws.Cells.Item(29, 3).Value = "4"
$ws.Cells.Item(29, 4).Value = "Quantum"
$ws.Cells.Item(30, 1).Value = "march-7th"
$ws.Cells.Item(30, 2).Value = "Preservation"
$ws.Cells.Item(30, 3).Value = "4"
$ws.Cells.Item(30, 4).Value = "Ice"
$ws.Cells.Item(31, 1).Value = "misha"
$ws.Cells.Item(31, 2).Value = "Destruction"
$ws.Cells.Item(31, 3).Value = "4"
$ws.Cells.Item(31, 4).Value = "Ice"
$ws.Cells.Item(32, 1).Value = "natasha"
$ws.Cells.Item(32, 2).Value = "Abundance"
$ws.Cells.Item(32, 3).Value = "4"
$ws.Cells.Item(32, 4).Value = "Physical"
$ws.Cells.Item(33, 1).Value = "pela"
$ws.Cells.Item(33, 2).Value = "Nihility"
$ws.Cells.Item(33, 3).Value = "4"
$ws.Cells.Item(33, 4).Value = "Ice"
$ws.Cells.Item(34, 1).Value = "qingque"
$ws.Cells.Item(34, 2).Value = "Erudition"
$ws.Cells.Item(34, 3).Value = "4"
$ws.Cells.Item(34, 4).Value = "Quantum"
$ws.Cells.Item(35, 1).Value = "ruan-mei"
$ws.Cells.Item(35, 2).Value = "Harmony"
$ws.Cells.Item(35, 3).Value = "5"
$ws.Cells.Item(35, 4).Value = "Ice"
$ws.Cells.Item(36, 1).Value = "sampo"
$ws.Cells.Item(36, 2).Value = "Nihility"
$ws.Cells.Item(36, 3).Value = "4"
$ws.Cells.Item(36, 4).Value = "Wind"
$ws.Cells.Item(37, 1).Value = "seele"
$ws.Cells.Item(37, 2).Value = "Hunt"
$ws.Cells.Item(37, 3).Value = "5"
$ws.Cells.Item(37, 4).Value = "Quantum"
$ws.Cells.Item(38, 1).Value = "serval"
$ws.Cells.Item(38, 2).Value = "Erudition"
$ws.Cells.Item(38, 3).Value = "4"
$ws.Cells.Item(38, 4).Value = "Lightning"
$ws.Cells.Item(39, 1).Value = "silver-wolf"
$ws.Cells.Item(39, 2).Value = "Nihility"
$ws.Cells.Item(39, 3).Value = "5"
$ws.Cells.Item(39, 4).Value = "Quantum"
$ws.Cells.Item(40, 1).Value = "sparkle"
$ws.Cells.Item(40, 2).Value = "Harmony"
$ws.Cells.Item(40, 3).Value = "5"
$ws.Cells.Item(40, 4).Value = "Quantum"
$ws.Cells.Item(41, 1).Value = "sushang"
$ws.Cells.Item(41, 2).Value = "Hunt"
$ws.Cells.Item(41, 3).Value = "4"
$ws.Cells.Item(41, 4).Value = "Physical"
$ws.Cells.Item(42, 1).Value = "tingyun"
$ws.Cells.Item(42, 2).Value = "Harmony"
$ws.Cells.Item(42, 3).Value = "4"
$ws.Cells.Item(42, 4).Value = "Lightning"
$ws.Cells.Item(43, 1).Value = "topaz"
$ws.Cells.Item(43, 2).Value = "Hunt"
$ws.Cells.Item(43, 3).Value = "5"
$ws.Cells.Item(43, 4).Value = "Fire"
$ws.Cells.Item(44, 1).Value = "trailblazer-fire"
$ws.Cells.Item(44, 2).Value = "Preservation"
$ws.Cells.Item(44, 3).Value = "5"
$ws.Cells.Item(44, 4).Value = "Fire"
$ws.Cells.Item(45, 1).Value = "trailblazer-physical"
$ws.Cells.Item(45, 2).Value = "Destruction"
$ws.Cells.Item(45, 3).Value = "5"
$ws.Cells.Item(45, 4).Value = "Physical"
$ws.Cells.Item(46, 1).Value = "welt"
$ws.Cells.Item(46, 2).Value = "Nihility"
$ws.Cells.Item(46, 3).Value = "5"
$ws.Cells.Item(46, 4).Value = "Imaginary"
$ws.Cells.Item(47, 1).Value = "xueyi"
$ws.Cells.Item(47, 2).Value = "Destruction"
$ws.Cells.Item(47, 3).Value = "4"
$ws.Cells.Item(47, 4).Value = "Quantum"
$ws.Cells.Item(48, 1).Value = "yanqing"
$ws.Cells.Item(48, 2).Value = "Hunt"
$ws.Cells.Item(48, 3).Value = "5"
$ws.Cells.Item(48, 4).Value = "Ice"
$ws.Cells.Item(49, 1).Value = "yukong"
$ws.Cells.Item(49, 2).Value = "Harmony"
$ws.Cells.Item(49, 3).Value = "4"
$ws.Cells.Item(49, 4).Value = "Imaginary"
